$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H18").Value = 125002470
$ws.Range("J18").Value = 4222
$ws.Range("L18").Value = 4222
$ws.Range("N18").Value = -4790
$ws.Range("H51").Value = 8214.556
$ws.Range("I51").Value = 3995.6667
$ws.Range("J51").Value = 10324
$ws.Range("K51").Value = 3995.6667
$ws.Range("L51").Value = 10324
$ws.Range("M51").Value = -3511.6667
$ws.Range("N51").Value = -11292
$ws.Range("H53").Value = 585.9091
$ws.Range("I53").Value = 1265
$ws.Range("K53").Value = 1265
$ws.Range("M53").Value = -628
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H98").Value = 912.5
$ws.Range("I98").Value = 933.3333
$ws.Range("K98").Value = 933.3333
$ws.Range("M98").Value = 564.6667
$ws.Range("H113").Value = 4000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10508
$ws.Range("H116").Value = 13312.375
$ws.Range("I116").Value = 10499.75
$ws.Range("J116").Value = 16125
$ws.Range("K116").Value = 10499.75
$ws.Range("L116").Value = 16125
$ws.Range("M116").Value = -7057.75
$ws.Range("N116").Value = -23009
$ws.Range("H122").Value = 912.5
$ws.Range("I122").Value = 933.3333
$ws.Range("K122").Value = 2799.9999
$ws.Range("M122").Value = -349.9998999999998
$ws.Range("H132").Value = 8717.031999999999
$ws.Range("I132").Value = 3236.0833
$ws.Range("K132").Value = 9708.249899999999
$ws.Range("M132").Value = -7178.249899999999
$ws.Range("H136").Value = 179992.28
$ws.Range("J136").Value = 179992.28
$ws.Range("L136").Value = 179992.28
$ws.Range("N136").Value = -190192.28
$ws.Range("H138").Value = 2199.3784
$ws.Range("I138").Value = 2090.4614
$ws.Range("K138").Value = 6271.3842
$ws.Range("M138").Value = -1131.3842

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 663.4375
$ws.Range("I4").Value = 431.14285
$ws.Range("J4").Value = 844.1111
$ws.Range("K4").Value = 431.14285
$ws.Range("L4").Value = 844.1111
$ws.Range("M4").Value = -315.14285
$ws.Range("N4").Value = -1076.1111
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H122").Value = 3954.9285
$ws.Range("J122").Value = 6799.6665
$ws.Range("L122").Value = 20398.9995
$ws.Range("N122").Value = -25298.9995

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4093.7856
$ws.Range("I134").Value = 4130.3413
$ws.Range("J134").Value = 2595
$ws.Range("K134").Value = 12391.0239
$ws.Range("L134").Value = 7785
$ws.Range("M134").Value = -9856.0239
$ws.Range("N134").Value = -12855

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 7500
$ws.Range("I13").Value = 3000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = -2861
$ws.Range("H29").Value = 11235.333
$ws.Range("I29").Value = 11189.571
$ws.Range("K29").Value = 11189.571
$ws.Range("M29").Value = -10896.571
$ws.Range("H31").Value = 4180.722
$ws.Range("I31").Value = 3473.8333
$ws.Range("J31").Value = 4534.1665
$ws.Range("K31").Value = 3473.8333
$ws.Range("L31").Value = 4534.1665
$ws.Range("M31").Value = -3178.8333
$ws.Range("N31").Value = -5124.1665
$ws.Range("H34").Value = 4180.722
$ws.Range("I34").Value = 3473.8333
$ws.Range("J34").Value = 4534.1665
$ws.Range("K34").Value = 3473.8333
$ws.Range("L34").Value = 4534.1665
$ws.Range("M34").Value = -3271.8333
$ws.Range("N34").Value = -4938.1665
$ws.Range("H99").Value = 6867.095
$ws.Range("J99").Value = 11936.875
$ws.Range("L99").Value = 11936.875
$ws.Range("N99").Value = -14932.875
$ws.Range("H126").Value = 6867.095
$ws.Range("J126").Value = 11936.875
$ws.Range("L126").Value = 35810.625
$ws.Range("N126").Value = -40750.625
$ws.Range("H134").Value = 15160.611
$ws.Range("I134").Value = 17657.584
$ws.Range("K134").Value = 52972.75199999999
$ws.Range("M134").Value = -50437.75199999999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1125.1578
$ws.Range("I11").Value = 688.2857
$ws.Range("J11").Value = 2348.4
$ws.Range("K11").Value = 2064.8571
$ws.Range("L11").Value = 7045.200000000001
$ws.Range("M11").Value = -1924.8571
$ws.Range("N11").Value = -7325.200000000001
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H121").Value = 1298.8572
$ws.Range("I121").Value = 25
$ws.Range("J121").Value = 1396.8462
$ws.Range("K121").Value = 75
$ws.Range("L121").Value = 4190.5386
$ws.Range("M121").Value = 1235
$ws.Range("N121").Value = -6810.5386
$ws.Range("H137").Value = 12560.737
$ws.Range("I137").Value = 9932.4
$ws.Range("K137").Value = 29797.2
$ws.Range("M137").Value = -24697.2

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 32000
$ws.Range("J48").Value = 32000
$ws.Range("L48").Value = 32000
$ws.Range("N48").Value = -32970
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 5135.577
$ws.Range("I132").Value = 4659.1177
$ws.Range("J132").Value = 6035.5557
$ws.Range("K132").Value = 13977.3531
$ws.Range("L132").Value = 18106.6671
$ws.Range("M132").Value = -11447.3531
$ws.Range("N132").Value = -23166.6671
$ws.Range("H133").Value = 44744.5
$ws.Range("I133").Value = 19709
$ws.Range("J133").Value = 69780
$ws.Range("K133").Value = 19709
$ws.Range("L133").Value = 69780
$ws.Range("M133").Value = -14649
$ws.Range("N133").Value = -79900

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 966
$ws.Range("I22").Value = 972.3333
$ws.Range("K22").Value = 972.3333
$ws.Range("M22").Value = -677.3333
$ws.Range("H27").Value = 966
$ws.Range("I27").Value = 972.3333
$ws.Range("K27").Value = 972.3333
$ws.Range("M27").Value = -865.3333
$ws.Range("H108").Value = 49999
$ws.Range("J108").Value = 49999
$ws.Range("L108").Value = 49999
$ws.Range("N108").Value = -57679

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 33999.5
$ws.Range("J47").Value = 33999.5
$ws.Range("L47").Value = 33999.5
$ws.Range("N47").Value = -35143.5
$ws.Range("H52").Value = 17499.5
$ws.Range("J52").Value = 29999
$ws.Range("L52").Value = 29999
$ws.Range("N52").Value = -30451
$ws.Range("H107").Value = 1155.5555
$ws.Range("I107").Value = 1442.3334
$ws.Range("K107").Value = 4327.0002
$ws.Range("M107").Value = -2407.0002
$ws.Range("H122").Value = 10964.565
$ws.Range("I122").Value = 7476.1763
$ws.Range("K122").Value = 22428.5289
$ws.Range("M122").Value = -19978.5289
$ws.Range("H132").Value = 3979.6775
$ws.Range("I132").Value = 3979.6775
$ws.Range("K132").Value = 11939.0325
$ws.Range("M132").Value = -9409.032499999999
